{"js": "// Remove the \"Overall detection probability for coyotes ...\" paragraph\n// that immediately follows the \"Estimated detection probabilities\" heading.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst target = paras.items.find(p =>\n  p.text && p.text.indexOf(\"Overall detection probability for coyotes\") !== -1\n);\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the \"Overall detection probability for coyotes ...\" paragraph\n# that immediately follows the \"Estimated detection probabilities\" heading.\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Overall detection probability for coyotes*\") {\n        $p.Range.Delete()\n    }\n}\n"}
